$d = $word.ActiveDocument
$nl = [char]11  # Word manual line break character -> serializes as <w:br/>

# "Objetivos" (PT/EN) paragraphs receive the text that used to live in
# "Programa resumido" (PT/EN); the reverse move happens further below.
$d.Paragraphs.Item(6).Range.Text = 'Processo Químico e Indústria Química; Química Fina; Petróleo, Gás Natural e Petroquímica; Plásticos e afins; Fertilizantes; Vidro; Celulose e Papel.'
$d.Paragraphs.Item(7).Range.Text = 'Chemical Process and Chemical Industry; Fine Chemistry; Oil, Natural Gas and Petrochemicals; Plastics and allied products; fertilizers; Glass; Cellulose and paper.'

# The bullet under "Docente(s) Responsavel(eis)" now holds the original
# "Objetivos" (PT) paragraph text.
$d.Paragraphs.Item(9).Range.Text = 'Conferir aos alunos uma visão geral da indústria química e correlatas, bem como das principais características dos processos desta indústria.'

# "Programa resumido" (PT/EN) paragraphs receive the original "Objetivos" (PT/EN) text.
$d.Paragraphs.Item(11).Range.Text = 'O conteúdo desta disciplina será de acordo com os tópicos a serem programados, devendo abordar assuntos relevantes relacionados a processos químicos e correlatas.'
$d.Paragraphs.Item(12).Range.Text = 'Providing to the students an overview of the chemical and related industries, as well as the main features of the processes and production arrangements of this industry.'

# "Programa" paragraph receives the original "Metodo" bullet body text.
$d.Paragraphs.Item(14).Range.Text = 'Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos.'

# "Avaliacao" bullet list: rotate the Metodo/Criterio/Norma body runs.
# Processed back-to-front (Norma, then Criterio, then Metodo) so each Find
# target is still unique/untouched at the time it is searched for.
$p17 = $d.Paragraphs.Item(17)

# 1) "Norma de recuperacao:" body is replaced by the full Bibliografia block
#    (moved up from the paragraph that used to follow "Bibliografia").
$biblio = 'Ullmann’s encyclopedia of industrial chemistry; Editorial advisory board, Giuseppe Bellussi et al.; 7th, completely revised edition; Weinheim ; New York : WileyVCH, 2011.' + $nl + $nl + 'Encyclopedia of Chemical Processing; Edited by Sunggyu Lee; New York : Taylor & Francis, 2006.' + $nl + $nl + 'Kirk, Raymond Eller. Encyclopedia of chemical technology / Herman F.Mark et al. New York: John Wiley, 1984.' + $nl + $nl + 'Manual Econômico da Indústria Química - MEIQ / Centro de Pesquisas e Desenvolvimento; 8ed; Camaçari: CEPED, 2007.' + $nl + $nl + 'Shreve, R. Norris; BRINK JR., J. A. Indústrias de processos químicos. Tradução de Horácio Macedo; 4.ed. Rio de Janeiro: Editora Guanabara Koogan, 2008, c1997.' + $nl + $nl + 'Revistas:' + $nl + $nl + 'Brazilian Journal of Chemical Engineering, São Paulo, SP: Brazilian Society of Chemical Engineering, v. 11, n. 1, 1995-;'
$null = $p17.Range.Find.Execute('Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplicação de prova escrita de recuperação valendo 10,00 pontos.', $false, $false, $false, $false, $false, $true, 1, $false, $biblio, 2)

# 2) "Criterio:" body becomes the old "Norma de recuperacao:" body, gaining a trailing break.
$null = $p17.Range.Find.Execute('A nota (NOTA) será composta por uma destas opções: prova em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula. A estas opções será incorporado, para cada aluno, seu respectivo percentual de frequência no cálculo da nota final (NF), conforme a fórmula explicitada abaixo:' + $nl + 'NF = NOTA x % FREQ.' + $nl, $false, $false, $false, $false, $false, $true, 1, $false, 'Frequência mínima de 70% e nota igual ou superior a 3,00 e inferior a 5,00 possibilita aplicação de prova escrita de recuperação valendo 10,00 pontos.' + $nl, 2)

# 3) "Metodo:" body becomes the old "Criterio:" body.
$null = $p17.Range.Find.Execute('Aulas expositivas, desenvolvimento de exercícios em sala e fora de sala de aula, discussão de casos práticos.' + $nl, $false, $false, $false, $false, $false, $true, 1, $false, 'A nota (NOTA) será composta por uma destas opções: prova em sala, apresentações em sala, entrega de exercícios ou casos práticos elaborados fora de sala de aula. A estas opções será incorporado, para cada aluno, seu respectivo percentual de frequência no cálculo da nota final (NF), conforme a fórmula explicitada abaixo:' + $nl + 'NF = NOTA x % FREQ.' + $nl, 2)

# The former Bibliografia text paragraph now holds the relocated Docente(s) entry.
$d.Paragraphs.Item(19).Range.Text = '1285870 - Marcos Villela Barcza'

Write-Host "Done."
